$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header layout: Sample | Cell_type | Condition | Mouse | Passage
# (write B:E before A so new shared-string entries land in the same
#  alphabetical order used by the authoring workbook: Cell_type, Condition,
#  Mouse, Passage, Sample)
$ws.Range("B1").Value = "Cell_type"
$ws.Range("C1").Value = "Condition"
$ws.Range("D1").Value = "Mouse"
$ws.Range("E1").Value = "Passage"
$ws.Range("A1").Value = "Sample"

# Return the cursor to the top-left cell (matches the saved view state of
# the authored workbook, replacing the stale G25 selection).
$ws.Range("A1").Select() | Out-Null
